$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1145.5714
$ws.Range("I2").Value = 255.75
$ws.Range("J2").Value = 2332
$ws.Range("K2").Value = 255.75
$ws.Range("L2").Value = 2332
$ws.Range("M2").Value = -142.75
$ws.Range("N2").Value = -2558
$ws.Range("H12").Value = 235.44444
$ws.Range("I12").Value = 202.125
$ws.Range("K12").Value = 202.125
$ws.Range("M12").Value = -32.125
$ws.Range("H19").Value = 1912.25
$ws.Range("I19").Value = 1824.5
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 1824.5
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -1649.5
$ws.Range("N19").Value = -2350
$ws.Range("H54").Value = 211479.8
$ws.Range("I54").Value = 253361
$ws.Range("K54").Value = 253361
$ws.Range("M54").Value = -252875
$ws.Range("H58").Value = 5183
$ws.Range("I58").Value = 3978.75
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 11936.25
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -11786.25
$ws.Range("N58").Value = -30300
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("N68").Value = -41498
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("N71").Value = -127488
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H111").Value = 1002
$ws.Range("I111").Value = 952.7143
$ws.Range("J111").Value = 1174.5
$ws.Range("K111").Value = 2858.1429
$ws.Range("L111").Value = 3523.5
$ws.Range("M111").Value = 208.8571000000002
$ws.Range("N111").Value = -9657.5
$ws.Range("H112").Value = 2077.9546
$ws.Range("J112").Value = 2105.476
$ws.Range("L112").Value = 6316.428
$ws.Range("N112").Value = -8532.428
$ws.Range("H128").Value = 37619
$ws.Range("J128").Value = 32968
$ws.Range("L128").Value = 32968
$ws.Range("N128").Value = -42928

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5195.7964
$ws.Range("I32").Value = 3651.48
$ws.Range("K32").Value = 3651.48
$ws.Range("M32").Value = -3364.48
$ws.Range("H45").Value = 2955.4736
$ws.Range("I45").Value = 2508.75
$ws.Range("K45").Value = 2508.75
$ws.Range("M45").Value = -2131.75
$ws.Range("H74").Value = 55561464
$ws.Range("I74").Value = 83336440
$ws.Range("J74").Value = 11514
$ws.Range("K74").Value = 83336440
$ws.Range("L74").Value = 11514
$ws.Range("M74").Value = -83335566
$ws.Range("N74").Value = -13262
$ws.Range("H77").Value = 55561464
$ws.Range("I77").Value = 83336440
$ws.Range("J77").Value = 11514
$ws.Range("K77").Value = 416682200
$ws.Range("L77").Value = 57570
$ws.Range("M77").Value = -416677832
$ws.Range("N77").Value = -66306

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2799.889
$ws.Range("I99").Value = 1457
$ws.Range("K99").Value = 1457
$ws.Range("M99").Value = 41
$ws.Range("H107").Value = 619.2857
$ws.Range("I107").Value = 539.7895
$ws.Range("J107").Value = 1374.5
$ws.Range("K107").Value = 539.7895
$ws.Range("L107").Value = 1374.5
$ws.Range("M107").Value = 1380.2105
$ws.Range("N107").Value = -5214.5
$ws.Range("H134").Value = 2328.1282
$ws.Range("I134").Value = 1436.1538
$ws.Range("J134").Value = 4112.077
$ws.Range("K134").Value = 4308.4614
$ws.Range("L134").Value = 12336.231
$ws.Range("M134").Value = -1773.4614
$ws.Range("N134").Value = -17406.231

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2469.8333
$ws.Range("I5").Value = 3535
$ws.Range("K5").Value = 3535
$ws.Range("M5").Value = -3423
$ws.Range("H8").Value = 1750
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("N8").Value = -1280
$ws.Range("H22").Value = 8671.166999999999
$ws.Range("I22").Value = 4009
$ws.Range("K22").Value = 4009
$ws.Range("M22").Value = -3659
$ws.Range("H25").Value = 5001
$ws.Range("I25").Value = 5001
$ws.Range("K25").Value = 5001
$ws.Range("M25").Value = -4827
$ws.Range("H86").Value = 4515.4707
$ws.Range("I86").Value = 4171.1816
$ws.Range("J86").Value = 5146.6665
$ws.Range("K86").Value = 4171.1816
$ws.Range("L86").Value = 5146.6665
$ws.Range("M86").Value = -3048.1816
$ws.Range("N86").Value = -7392.6665
$ws.Range("H89").Value = 4515.4707
$ws.Range("I89").Value = 4171.1816
$ws.Range("J89").Value = 5146.6665
$ws.Range("K89").Value = 20855.908
$ws.Range("L89").Value = 25733.3325
$ws.Range("M89").Value = -15239.908
$ws.Range("N89").Value = -36965.3325

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 343751
$ws.Range("J37").Value = 343751
$ws.Range("L37").Value = 1031253
$ws.Range("N37").Value = -1031477
$ws.Range("H113").Value = 76924616
$ws.Range("I113").Value = 1392.2858
$ws.Range("J113").Value = 166668370
$ws.Range("K113").Value = 4176.857400000001
$ws.Range("L113").Value = 500005110
$ws.Range("M113").Value = -2006.857400000001
$ws.Range("N113").Value = -500009450
$ws.Range("H122").Value = 8427091
$ws.Range("J122").Value = 8405825
$ws.Range("L122").Value = 75652425
$ws.Range("N122").Value = -75657325

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 1100000
$ws.Range("J47").Value = 1100000
$ws.Range("L47").Value = 1100000
$ws.Range("N47").Value = -1101136
$ws.Range("H70").Value = 10497.923
$ws.Range("J70").Value = 11498
$ws.Range("L70").Value = 11498
$ws.Range("N70").Value = -12038
$ws.Range("H73").Value = 10497.923
$ws.Range("J73").Value = 11498
$ws.Range("L73").Value = 11498
$ws.Range("N73").Value = -13370
$ws.Range("H102").Value = 1990.6857
$ws.Range("I102").Value = 1100.3334
$ws.Range("K102").Value = 1100.3334
$ws.Range("M102").Value = 521.6666
$ws.Range("H113").Value = 3202.2354
$ws.Range("I113").Value = 2293.8
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 2293.8
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -123.8000000000002
$ws.Range("N113").Value = -8840
$ws.Range("H122").Value = 11689
$ws.Range("I122").Value = 11215.5
$ws.Range("J122").Value = 15003.5
$ws.Range("K122").Value = 33646.5
$ws.Range("L122").Value = 45010.5
$ws.Range("M122").Value = -31196.5
$ws.Range("N122").Value = -49910.5
$ws.Range("H123").Value = 36885.2
$ws.Range("J123").Value = 36885.2
$ws.Range("L123").Value = 36885.2
$ws.Range("N123").Value = -41785.2
$ws.Range("H128").Value = 59666.668
$ws.Range("J128").Value = 59666.668
$ws.Range("L128").Value = 59666.668
$ws.Range("N128").Value = -69626.66800000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 21000
$ws.Range("J59").Value = 21000
$ws.Range("L59").Value = 21000
$ws.Range("N59").Value = -22308
$ws.Range("H61").Value = 6790.4614
$ws.Range("I61").Value = 5127.2
$ws.Range("K61").Value = 5127.2
$ws.Range("M61").Value = -4925.2
$ws.Range("H93").Value = 3531.2222
$ws.Range("I93").Value = 3580.875
$ws.Range("K93").Value = 3580.875
$ws.Range("M93").Value = -2332.875
$ws.Range("H113").Value = 6790.4614
$ws.Range("I113").Value = 5127.2
$ws.Range("K113").Value = 5127.2
$ws.Range("M113").Value = -2957.2

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 271.42856
$ws.Range("I113").Value = 209.23529
$ws.Range("J113").Value = 535.75
$ws.Range("K113").Value = 627.70587
$ws.Range("L113").Value = 1607.25
$ws.Range("M113").Value = 1542.29413
$ws.Range("N113").Value = -5947.25
$ws.Range("H119").Value = 78947
$ws.Range("J119").Value = 78947
$ws.Range("L119").Value = 78947
$ws.Range("N119").Value = -88623
$ws.Range("H132").Value = 2260.9167
$ws.Range("I132").Value = 1696.4445
$ws.Range("K132").Value = 5089.333500000001
$ws.Range("M132").Value = -2559.333500000001
$ws.Range("H137").Value = 69684.28999999999
$ws.Range("J137").Value = 69684.28999999999
$ws.Range("L137").Value = 69684.28999999999
$ws.Range("N137").Value = -79884.28999999999
